$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update price and volume columns for rows with straightforward value changes
$ws.Range("D2").Value = "43.768.13"
$ws.Range("E2").Value = "  +0.19%  "
$ws.Range("D3").Value = "2.223.94"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.14%  "
$ws.Range("D5").Value = "'271.04"
$ws.Range("E5").Value = "  +4.78%  "
$ws.Range("D6").Value = "'93.60"
$ws.Range("E6").Value = "  +15.39%  "
$ws.Range("D7").Value = "'0.628"
$ws.Range("E7").Value = "  +0.50%  "
$ws.Range("E8").Value = "  +0.06%  "
$ws.Range("D9").Value = "'0.624"
$ws.Range("E9").Value = "  +4.84%  "
$ws.Range("D10").Value = "'46.04"
$ws.Range("E10").Value = "  +6.72%  "
$ws.Range("D11").Value = "'0.0979"
$ws.Range("E11").Value = "  +6.35%  "
$ws.Range("D12").Value = "'8.29"
$ws.Range("E12").Value = "  +19.37%  "
$ws.Range("E13").Value = "  +1.63%  "
$ws.Range("D18").Value = "43.729.33"
$ws.Range("E18").Value = "  +0.29%  "
$ws.Range("E19").Value = "  +2.74%  "
$ws.Range("E20").Value = "  +1.83%  "
$ws.Range("D21").Value = "'70.67"
$ws.Range("E21").Value = "  +0.48%  "
$ws.Range("E22").Value = "  -1.14%  "
$ws.Range("D23").Value = "'233.87"
$ws.Range("E23").Value = "  +1.46%  "
$ws.Range("D24").Value = "'9.15"
$ws.Range("E24").Value = "  +2.40%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'11.44"
$ws.Range("E26").Value = "  +7.47%  "
$ws.Range("D27").Value = "'2.51"
$ws.Range("E27").Value = "  +12.75%  "
$ws.Range("D28").Value = "'41.32"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("D29").Value = "'3.52"
$ws.Range("E29").Value = "  +4.83%  "
$ws.Range("E30").Value = "  +2.00%  "
$ws.Range("D31").Value = "'172.57"
$ws.Range("E31").Value = "  -0.25%  "
$ws.Range("E32").Value = "  +5.21%  "
$ws.Range("D33").Value = "'20.93"
$ws.Range("E33").Value = "  +2.67%  "
$ws.Range("E34").Value = "  +4.49%  "
$ws.Range("E35").Value = "  +1.55%  "
$ws.Range("E36").Value = "  -1.90%  "
$ws.Range("E37").Value = "  -0.48%  "
$ws.Range("E38").Value = "  -3.34%  "
$ws.Range("D39").Value = "'3.57"
$ws.Range("E39").Value = "  +24.30%  "
$ws.Range("D40").Value = "'12.91"
$ws.Range("E40").Value = "  -2.17%  "
$ws.Range("D41").Value = "'0.223"
$ws.Range("E41").Value = "  +11.93%  "
$ws.Range("E42").Value = "  +2.66%  "
$ws.Range("D43").Value = "'63.78"
$ws.Range("E43").Value = "  +2.19%  "
$ws.Range("D44").Value = "'5.34"
$ws.Range("E44").Value = "  -1.92%  "
$ws.Range("E45").Value = "  +0.62%  "
$ws.Range("D46").Value = "'8.36"
$ws.Range("E46").Value = "  +1.45%  "
$ws.Range("D47").Value = "'100.43"
$ws.Range("E47").Value = "  -0.63%  "
$ws.Range("E48").Value = "  +4.45%  "
$ws.Range("E49").Value = "  +1.94%  "
$ws.Range("D50").Value = "'0.442"
$ws.Range("E50").Value = "  +0.74%  "
$ws.Range("D51").Value = "'1.45"
$ws.Range("E51").Value = "  -5.04%  "

# Rows 14-15 swap rank order between Chainlink and WrappedliquidstakedEther2.0, with updated prices
$ws.Range("B14").Value = "Chainlink"
$ws.Range("C14").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D14").Value = "'15.09"
$ws.Range("E14").Value = "  +5.53%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.554.16"
$ws.Range("E15").Value = "  +1.55%  "

# Rows 16-17 swap rank order between Polygon and WrappedEther, with updated prices
$ws.Range("B16").Value = "Polygon"
$ws.Range("C16").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D16").Value = "'0.807"
$ws.Range("E16").Value = "  +4.02%  "
$ws.Range("B17").Value = "WrappedEther"
$ws.Range("C17").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D17").Value = "2.208.97"
$ws.Range("E17").Value = "  +0.45%  "
